# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row -> new F value
$sheet1Updates = @{
    3  = 593
    7  = 14827
    8  = 403
    10 = 15199
    11 = 33
    12 = 8692
    13 = 329
    15 = 63
    16 = 180
    20 = 15
    21 = 4
    22 = 20
    26 = 7
    27 = 12
    28 = 55
    32 = 27
    35 = 266
    36 = 422
    38 = 5346
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Sheet "全部类型" (sheet4): row -> new F value
$sheet4Updates = @{
    3  = 593
    7  = 14827
    8  = 403
    10 = 15199
    11 = 33
    12 = 8692
    13 = 329
    16 = 63
    17 = 180
    21 = 15
    22 = 4
    23 = 20
    27 = 7
    28 = 12
    29 = 55
    35 = 27
    38 = 266
    39 = 422
    41 = 5346
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
